# Applies the "DLC new skills have being supported" update to the gear
# calculation workbook: new numeric inputs for existing fields plus two new
# armor-skill columns (巧击 / 偷袭) and a new 弱特属性 skill row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Plain numeric value updates on already-existing (already styled) cells
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 330

$ws.Range("B6").Value = 0.63
$ws.Range("C6").Value = 0.37

$ws.Range("A8").Value = 80

$ws.Range("A10").Value = 57

$ws.Range("A12").Value = 1.3

$ws.Range("A16").Value = 51
$ws.Range("B16").Value = 25

# Updated instructional text (readme -> readme.md)
$ws.Range("A18").Value = "防具技能选择，左边填最低允许等级，右边填最高允许等级，不需要此技能请至少将左边填“-1”  覆盖率为0-1的小数，请注意相关技能之间覆盖率的彼此影响，具体请查看readme.md文档"

# Row 20 (skill level selections)
$ws.Range("A20").Value = 4
$ws.Range("C20").Value = -1
$ws.Range("E20").Value = 1
$ws.Range("G20").Value = -1
$ws.Range("I20").Value = 2
$ws.Range("K20").Value = -1
$ws.Range("M20").Value = -1
$ws.Range("O20").Value = 1
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = -1
$ws.Range("R20").Value = 2
$ws.Range("U20").Value = -1
$ws.Range("V20").Value = 2

# Row 22 (coverage rates)
$ws.Range("C22").Value = 0.4
$ws.Range("I22").Value = 0.76
$ws.Range("O22").Value = 0.6
$ws.Range("Q22").Value = 0.3
$ws.Range("W22").Value = 0.5

# Row 28
$ws.Range("A28").Value = 4
$ws.Range("C28").Value = -1
$ws.Range("E28").Value = -1

# Row 30
$ws.Range("C30").Value = 0.8

# ---------------------------------------------------------------------
# 2. New armor-skill columns Y:AB for 巧击 / 偷袭 (mirrors existing layout)
# ---------------------------------------------------------------------

# Headers (row 19) - plain text cells, default (unfilled) style
$ws.Range("Y19").Value = "巧击"
$ws.Range("Z19").Value = "（0，3）"
$ws.Range("AA19").Value = "偷袭"
$ws.Range("AB19").Value = "（0，3）"

# Row 21 coverage labels - plain text cells
$ws.Range("Y21").Value = "覆盖率"
$ws.Range("AA21").Value = "覆盖率"

# ---------------------------------------------------------------------
# 3. New 弱特属性 skill (row 31/33 labels, row 32/34 data) under E:F
# ---------------------------------------------------------------------
$ws.Range("E31").Value = "弱特属性"
$ws.Range("F31").Value = "（0，3）"

$ws.Range("E33").Value = "覆盖率"

# ---------------------------------------------------------------------
# 4. Fill-styled input cells: copy the format used by the other highlighted
#    input cells (e.g. A2) onto the newly created cells, then set values.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("Y20:AB20").PasteSpecial(-4122)
$ws.Range("Y22").PasteSpecial(-4122)
$ws.Range("AA22").PasteSpecial(-4122)
$ws.Range("E32:F32").PasteSpecial(-4122)
$ws.Range("E34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("Y20").Value = 2
$ws.Range("Z20").Value = 2
$ws.Range("AA20").Value = 1
$ws.Range("AB20").Value = 1

$ws.Range("Y22").Value = 0.74
$ws.Range("AA22").Value = 0.5

$ws.Range("E32").Value = -1
$ws.Range("F32").Value = 3

$ws.Range("E34").Value = 0.7

# ---------------------------------------------------------------------
# 5. Selection / view bookkeeping
# ---------------------------------------------------------------------
$ws.Range("J28").Select()
